# Add a new row 4 to Sheet1:
#   A4 = 9          (same numeric/bold/bordered style as the existing A2/A3 cells)
#   E4 = "ghv"
#   B4,C4,D4,F4,G4,H4,I4,J4,K4 stay blank (mirrors the mostly-empty row added upstream)
# This grows the sheet's used range from A1:K3 to A1:K4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric value for column A, row 4.
$ws.Range("A4").Value = 9

# Match the formatting already used by the other data rows' column-A cells
# (bold font, thin border, centered) by copying A3's format onto A4.
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122)   # xlPasteFormats

# Text value for column E, row 4.
$ws.Range("E4").Value = "ghv"

# Leave B4, C4, D4, F4, G4, H4, I4, J4, K4 empty - no content for these cells,
# same as the new row in the target workbook (only A4 and E4 carry data).
